$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the specific Price cells we are about to rewrite to be treated as
# text *before* assigning their new values. Several of these values look
# numeric (e.g. "1.00", "0.171", "66.796.87") and Excel would otherwise
# silently reinterpret them as floating point numbers (losing trailing
# zeros / introducing FP rounding), which would not match the original
# inline-string text content of these cells. (Looping per-cell because a
# single multi-area Range(...).NumberFormat assignment only affects the
# first area in this engine.)
$priceRows = @(2,3,5,6,7,8,9,10,11,12,13,14,15,16,17,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,37,38,39,40,41,42,43,44,45,46,47,48,49,51)
foreach ($r in $priceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = '66.796.87'
$ws.Range("E2").Value = '  -4.96%  '

$ws.Range("D3").Value = '3.223.55'
$ws.Range("E3").Value = '  -8.47%  '

$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").Value = '584.25'
$ws.Range("E5").Value = '  -3.86%  '

$ws.Range("D6").Value = '151.66'
$ws.Range("E6").Value = '  -12.44%  '

$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").Value = '3.212.88'
$ws.Range("E8").Value = '  -8.62%  '

$ws.Range("D9").Value = '0.547'
$ws.Range("E9").Value = '  -9.98%  '

$ws.Range("D10").Value = '0.171'
$ws.Range("E10").Value = '  -12.66%  '

$ws.Range("D11").Value = '6.48'
$ws.Range("E11").Value = '  -9.56%  '

$ws.Range("D12").Value = '0.503'
$ws.Range("E12").Value = '  -14.47%  '

$ws.Range("D13").Value = '38.92'
$ws.Range("E13").Value = '  -15.94%  '

$ws.Range("D14").Value = '0.0000244'
$ws.Range("E14").Value = '  -11.65%  '

$ws.Range("D15").Value = '3.733.49'
$ws.Range("E15").Value = '  -8.70%  '

$ws.Range("D16").Value = '66.818.03'
$ws.Range("E16").Value = '  -4.96%  '

$ws.Range("D17").Value = '3.225.41'
$ws.Range("E17").Value = '  -9.01%  '

$ws.Range("E18").Value = '  -5.56%  '

$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = '535.86'
$ws.Range("E19").Value = '  -12.60%  '

$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").Value = '7.24'
$ws.Range("E20").Value = '  -14.02%  '

$ws.Range("D21").Value = '15.22'
$ws.Range("E21").Value = '  -14.21%  '

$ws.Range("D22").Value = '0.764'
$ws.Range("E22").Value = '  -13.37%  '

$ws.Range("D23").Value = '7.86'
$ws.Range("E23").Value = '  -13.69%  '

$ws.Range("D24").Value = '13.69'
$ws.Range("E24").Value = '  -12.08%  '

$ws.Range("D25").Value = '85.43'
$ws.Range("E25").Value = '  -13.58%  '

$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  -0.21%  '

$ws.Range("D27").Value = '3.18'
$ws.Range("E27").Value = '  -15.31%  '

$ws.Range("D28").Value = '2.19'
$ws.Range("E28").Value = '  -15.14%  '

$ws.Range("D29").Value = '8.13'
$ws.Range("E29").Value = '  -10.72%  '

$ws.Range("D30").Value = '29.27'
$ws.Range("E30").Value = '  -13.18%  '

$ws.Range("D31").Value = '2.59'
$ws.Range("E31").Value = '  -13.81%  '

$ws.Range("D32").Value = '1.14'
$ws.Range("E32").Value = '  -12.54%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '6.54'
$ws.Range("E33").Value = '  -19.47%  '

$ws.Range("B34").Value = 'Bittensor'
$ws.Range("C34").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D34").Value = '534.66'
$ws.Range("E34").Value = '  -15.13%  '

$ws.Range("D35").Value = '5.75'
$ws.Range("E35").Value = '  -15.67%  '

$ws.Range("E36").Value = '  +0.30%  '

$ws.Range("D37").Value = '53.36'
$ws.Range("E37").Value = '  -5.82%  '

$ws.Range("D38").Value = '0.0862'
$ws.Range("E38").Value = '  -14.14%  '

$ws.Range("D39").Value = '0.0422'
$ws.Range("E39").Value = '  -17.05%  '

$ws.Range("D40").Value = '9.30'
$ws.Range("E40").Value = '  -13.72%  '

$ws.Range("D41").Value = '0.126'
$ws.Range("E41").Value = '  -14.98%  '

$ws.Range("D42").Value = '2.928.84'
$ws.Range("E42").Value = '  -13.39%  '

$ws.Range("D43").Value = '2.59'
$ws.Range("E43").Value = '  -25.16%  '

$ws.Range("D44").Value = '0.264'
$ws.Range("E44").Value = '  -15.24%  '

$ws.Range("D45").Value = '0.0₃0583'
$ws.Range("E45").Value = '  -22.13%  '

$ws.Range("D46").Value = '2.38'
$ws.Range("E46").Value = '  -17.36%  '

$ws.Range("B47").Value = 'USDe'
$ws.Range("C47").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D47").Value = '1.00'
$ws.Range("E47").Value = '  -0.08%  '

$ws.Range("D48").Value = '2.13'
$ws.Range("E48").Value = '  -16.60%  '

$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = '26.03'
$ws.Range("E49").Value = '  -18.92%  '

$ws.Range("E50").Value = '  -12.32%  '

$ws.Range("D51").Value = '122.75'
$ws.Range("E51").Value = '  -8.14%  '
